$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 609, shifting existing rows 609-680 down to 610-681.
$ws.Rows.Item(609).Insert()

# Populate the newly inserted row 609 with the new record's data.
$ws.Range("A609").Value = 3
$ws.Range("B609").Value = "Femacal de La Calera"
$ws.Range("C609").Value = "Coquimbo"
$ws.Range("D609").Value = 45124
$ws.Range("E609").Value = 5
$ws.Range("F609").Value = 100112017
$ws.Range("G609").Value = "Apio"
$ws.Range("H609").Value = "Americana (o)"
$ws.Range("I609").Value = "Primera"
$ws.Range("J609").Value = 110
$ws.Range("K609").Value = 8500
$ws.Range("L609").Value = 8500
$ws.Range("M609").Value = 8500
$ws.Range("N609").Value = "`$/docena de matas"
$ws.Range("O609").Value = "Provincia de Limarí"
$ws.Range("P609").Value = 1417
$ws.Range("Q609").Value = 6
$ws.Range("R609").Value = "Hortaliza"
